$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 51
$ws.Range("H51").Value = 2620.5454
$ws.Range("J51").Value = 3919.4
$ws.Range("L51").Value = 3919.4
$ws.Range("N51").Value = -4887.4

# Row 132
$ws.Range("H132").Value = 4954.7144
$ws.Range("I132").Value = 5384.7393
$ws.Range("K132").Value = 16154.2179
$ws.Range("M132").Value = -13624.2179

# Row 137
$ws.Range("H137").Value = 5380.375
$ws.Range("I137").Value = 11126.667
$ws.Range("K137").Value = 33380.001
$ws.Range("M137").Value = -30830.001

# Row 138
$ws.Range("H138").Value = 6065420
$ws.Range("I138").Value = 2381.6
$ws.Range("K138").Value = 7144.799999999999
$ws.Range("M138").Value = -2004.799999999999

$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 13167552
$ws.Range("J32").Value = 23373.273
$ws.Range("L32").Value = 23373.273
$ws.Range("N32").Value = -23947.273

# Row 61
$ws.Range("H61").Value = 28574750
$ws.Range("I61").Value = 50002224
$ws.Range("K61").Value = 50002224
$ws.Range("M61").Value = -50002012

# Row 74
$ws.Range("H74").Value = 28606218
$ws.Range("I74").Value = 38506110
$ws.Range("J74").Value = 6531.8887
$ws.Range("K74").Value = 38506110
$ws.Range("L74").Value = 6531.8887
$ws.Range("M74").Value = -38505236
$ws.Range("N74").Value = -8279.8887

# Row 77
$ws.Range("H77").Value = 28606218
$ws.Range("I77").Value = 38506110
$ws.Range("J77").Value = 6531.8887
$ws.Range("K77").Value = 192530550
$ws.Range("L77").Value = 32659.4435
$ws.Range("M77").Value = -192526182
$ws.Range("N77").Value = -41395.4435

# Row 102
$ws.Range("H102").Value = 2827.25
$ws.Range("I102").Value = 2827.25
$ws.Range("J102").Value = 0
$ws.Range("K102").Value = 2827.25
$ws.Range("L102").Value = 0
$ws.Range("M102").ClearContents()
$ws.Range("N102").Value = -1205.25

# Row 122
$ws.Range("H122").Value = 2589.4062
$ws.Range("I122").Value = 1469.2941
$ws.Range("J122").Value = 3858.8667
$ws.Range("K122").Value = 4407.8823
$ws.Range("L122").Value = 11576.6001
$ws.Range("M122").Value = -1957.8823
$ws.Range("N122").Value = -16476.6001

# Row 132
$ws.Range("H132").Value = 71432184
$ws.Range("I132").Value = 3893.923
$ws.Range("J132").Value = 1000000000
$ws.Range("K132").Value = 11681.769
$ws.Range("L132").Value = 3000000000
$ws.Range("M132").Value = -9151.769
$ws.Range("N132").Value = -3000005060

# Row 136
$ws.Range("H136").Value = 28574750
$ws.Range("I136").Value = 50002224
$ws.Range("K136").Value = 150006672
$ws.Range("M136").Value = -150004122

$ws = $wb.Worksheets.Item("BSM")
# Row 20
$ws.Range("H20").Value = 3775.5
$ws.Range("I20").Value = 4108.2856
$ws.Range("J20").Value = 2999
$ws.Range("K20").Value = 4108.2856
$ws.Range("L20").Value = 2999
$ws.Range("M20").Value = -3861.2856
$ws.Range("N20").Value = -3493

# Row 80
$ws.Range("H80").Value = 3003.6875
$ws.Range("J80").Value = 5051.875
$ws.Range("L80").Value = 5051.875
$ws.Range("N80").Value = -7047.875

# Row 83
$ws.Range("H83").Value = 3003.6875
$ws.Range("J83").Value = 5051.875
$ws.Range("L83").Value = 25259.375
$ws.Range("N83").Value = -35243.375

$ws = $wb.Worksheets.Item("CRP")
# Row 134
$ws.Range("H134").Value = 1480.7391
$ws.Range("I134").Value = 1161.6316
$ws.Range("J134").Value = 2996.5
$ws.Range("K134").Value = 3484.8948
$ws.Range("L134").Value = 8989.5
$ws.Range("M134").Value = -949.8948
$ws.Range("N134").Value = -14059.5

$ws = $wb.Worksheets.Item("CUL")
# Row 121
$ws.Range("H121").Value = 1075.4667
$ws.Range("J121").Value = 1153.5555
$ws.Range("L121").Value = 3460.6665
$ws.Range("N121").Value = -6080.666499999999

# Row 131
$ws.Range("H131").Value = 39028.035
$ws.Range("I131").Value = 70307.87
$ws.Range("J131").Value = 7748.2
$ws.Range("K131").Value = 210923.61
$ws.Range("L131").Value = 23244.6
$ws.Range("M131").Value = -205883.61
$ws.Range("N131").Value = -33324.6

# Row 132
$ws.Range("H132").Value = 2473313
$ws.Range("I132").Value = 2874.625
$ws.Range("K132").Value = 25871.625
$ws.Range("M132").Value = -23341.625

$ws = $wb.Worksheets.Item("GSM")
# Row 113
$ws.Range("H113").Value = 3101.6316
$ws.Range("J113").Value = 3948.6667
$ws.Range("L113").Value = 3948.6667
$ws.Range("N113").Value = -8288.6667

# Row 132
$ws.Range("H132").Value = 2419.375
$ws.Range("I132").Value = 1935.125
$ws.Range("J132").Value = 2903.625
$ws.Range("K132").Value = 5805.375
$ws.Range("L132").Value = 8710.875
$ws.Range("M132").Value = -3275.375
$ws.Range("N132").Value = -13770.875

$ws = $wb.Worksheets.Item("LTW")
# Row 122
$ws.Range("H122").Value = 4613.391
$ws.Range("I122").Value = 4141.9165
$ws.Range("J122").Value = 5127.727
$ws.Range("K122").Value = 12425.7495
$ws.Range("L122").Value = 15383.181
$ws.Range("M122").Value = -9975.749500000002
$ws.Range("N122").Value = -20283.181

# Row 132
$ws.Range("H132").Value = 105265000
$ws.Range("I132").Value = 1763.7693
$ws.Range("J132").Value = 333335330
$ws.Range("K132").Value = 5291.3079
$ws.Range("L132").Value = 1000005990
$ws.Range("M132").Value = -2761.3079
$ws.Range("N132").Value = -1000011050

# Row 8
$ws.Range("H8").Value = 1500
$ws.Range("J8").Value = 1500
$ws.Range("L8").Value = 1500
$ws.Range("N8").Value = -1780

$ws = $wb.Worksheets.Item("WVR")
# Row 26
$ws.Range("H26").Value = 10000
$ws.Range("I26").Value = 10000
$ws.Range("K26").Value = 10000
$ws.Range("M26").Value = -9707

# Row 29
$ws.Range("H29").Value = 6499.6665
$ws.Range("I29").Value = 6499.6665
$ws.Range("J29").Value = 0
$ws.Range("K29").Value = 6499.6665
$ws.Range("L29").Value = 0
$ws.Range("M29").ClearContents()
$ws.Range("N29").Value = -6209.6665

# Row 126
$ws.Range("H126").Value = 3494.0857
$ws.Range("I126").Value = 4277.76
$ws.Range("J126").Value = 1534.9
$ws.Range("K126").Value = 12833.28
$ws.Range("L126").Value = 4604.700000000001
$ws.Range("M126").Value = -10363.28
$ws.Range("N126").Value = -9544.700000000001

# Row 132
$ws.Range("H132").Value = 3588.5398
$ws.Range("I132").Value = 3626.276
$ws.Range("J132").Value = 3150.8
$ws.Range("K132").Value = 10878.828
$ws.Range("L132").Value = 9452.400000000001
$ws.Range("M132").Value = -8348.828
$ws.Range("N132").Value = -14512.4

# Row 136
$ws.Range("H136").Value = 1634.8298
$ws.Range("I136").Value = 1665.0975
$ws.Range("J136").Value = 1428
$ws.Range("K136").Value = 4995.2925
$ws.Range("L136").Value = 4284
$ws.Range("M136").Value = -2445.2925
$ws.Range("N136").Value = -9384
